# remove column from alcohol data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete column M (the "alcohol" data column being removed); column N
# shifts left to become the new column M.
$ws.Columns.Item(13).Delete()

# Move the selection to the new last-used column (M1), matching the
# post-edit cursor position recorded in the file.
$ws.Range("M1").Select()
